$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3915976666666667
$ws.Range("H2").Value = 1.174793
$ws.Range("I2").Value = 0.02606065131430495
$ws.Range("J2").Value = 0.02606065131430495
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.66643366666667
$ws.Range("N2").Value = 37.999301
$ws.Range("O2").Value = 0.1759291503241684
$ws.Range("P2").Value = 0.1759291503241684
$ws.Range("Q2").Value = 4.960145868854779
$ws.Range("R2").Value = 44.64131281969301
$ws.Range("S2").Value = 0.004584828242620093
$ws.Range("T2").Value = 0.004584828242620093
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3915976666666667
$ws.Range("H3").Value = 1.174793
$ws.Range("I3").Value = 0.02606065131430495
$ws.Range("J3").Value = 0.02606065131430495
$ws.Range("O3").Value = 0.5164516272884614
$ws.Range("P3").Value = 0.5164516272884614
$ws.Range("Q3").Value = 14.56083543197945
$ws.Range("R3").Value = 131.047518887815
$ws.Range("S3").Value = 0.01345906577946997
$ws.Range("T3").Value = 0.01345906577946997
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3915976666666667
$ws.Range("H4").Value = 1.174793
$ws.Range("I4").Value = 0.02606065131430495
$ws.Range("J4").Value = 0.02606065131430495
$ws.Range("M4").Value = 22.14777066666666
$ws.Range("N4").Value = 66.44331199999999
$ws.Range("O4").Value = 0.3076192223873702
$ws.Range("P4").Value = 0.3076192223873702
$ws.Range("Q4").Value = 8.67301531493511
$ws.Range("R4").Value = 78.057137834416
$ws.Range("S4").Value = 0.008016757292214887
$ws.Range("T4").Value = 0.008016757292214887
$ws.Range("I5").Value = 0.4187506438669658
$ws.Range("J5").Value = 0.4187506438669658
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.66643366666667
$ws.Range("N5").Value = 37.999301
$ws.Range("O5").Value = 0.1759291503241684
$ws.Range("P5").Value = 0.1759291503241684
$ws.Range("Q5").Value = 79.70116522440435
$ws.Range("R5").Value = 717.310487019639
$ws.Range("S5").Value = 0.07367044497321375
$ws.Range("T5").Value = 0.07367044497321375
$ws.Range("I6").Value = 0.4187506438669658
$ws.Range("J6").Value = 0.4187506438669658
$ws.Range("O6").Value = 0.5164516272884614
$ws.Range("P6").Value = 0.5164516272884614
$ws.Range("S6").Value = 0.2162644514531855
$ws.Range("T6").Value = 0.2162644514531855
$ws.Range("I7").Value = 0.4187506438669658
$ws.Range("J7").Value = 0.4187506438669658
$ws.Range("M7").Value = 22.14777066666666
$ws.Range("N7").Value = 66.44331199999999
$ws.Range("O7").Value = 0.3076192223873702
$ws.Range("P7").Value = 0.3076192223873702
$ws.Range("Q7").Value = 139.3607052868853
$ws.Range("R7").Value = 1254.246347581968
$ws.Range("S7").Value = 0.1288157474405666
$ws.Range("T7").Value = 0.1288157474405666
$ws.Range("G8").Value = 8.342485333333334
$ws.Range("H8").Value = 25.027456
$ws.Range("I8").Value = 0.5551887048187292
$ws.Range("J8").Value = 0.5551887048187292
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.66643366666667
$ws.Range("N8").Value = 37.999301
$ws.Range("O8").Value = 0.1759291503241684
$ws.Range("P8").Value = 0.1759291503241684
$ws.Range("Q8").Value = 105.6695370898062
$ws.Range("R8").Value = 951.0258338082562
$ws.Range("S8").Value = 0.09767387710833457
$ws.Range("T8").Value = 0.09767387710833457
$ws.Range("G9").Value = 8.342485333333334
$ws.Range("H9").Value = 25.027456
$ws.Range("I9").Value = 0.5551887048187292
$ws.Range("J9").Value = 0.5551887048187292
$ws.Range("O9").Value = 0.5164516272884614
$ws.Range("P9").Value = 0.5164516272884614
$ws.Range("Q9").Value = 310.1998974262756
$ws.Range("R9").Value = 2791.79907683648
$ws.Range("S9").Value = 0.286728110055806
$ws.Range("T9").Value = 0.286728110055806
$ws.Range("G10").Value = 8.342485333333334
$ws.Range("H10").Value = 25.027456
$ws.Range("I10").Value = 0.5551887048187292
$ws.Range("J10").Value = 0.5551887048187292
$ws.Range("M10").Value = 22.14777066666666
$ws.Range("N10").Value = 66.44331199999999
$ws.Range("O10").Value = 0.3076192223873702
$ws.Range("P10").Value = 0.3076192223873702
$ws.Range("Q10").Value = 184.7674519526969
$ws.Range("R10").Value = 1662.907067574272
$ws.Range("S10").Value = 0.1707867176545887
$ws.Range("T10").Value = 0.1707867176545887
